# Edit the Sheet1 header row + data row to match the updated template:
#  - rename a few header labels (typo fixes / renames)
#  - add a "nama_uuk" column label in A1 (replacing the old "id_uuk")
#  - clear the now-unused last header cell (AG1)
#  - change A2 from a numeric placeholder to the text "asd"
#  - remove the now-unused AG2 data cell
#  - widen column P and refresh the window scroll/selection position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate() | Out-Null

# --- Header row (row 1) label fixes --------------------------------------
$ws.Range("A1").Value = "nama_uuk"
$ws.Range("G1").Value = "kontrak_tanggal"
$ws.Range("H1").Value = "kontrak_nomor"
$ws.Range("P1").Value = "keuangan_sisa_invoice_total"
$ws.Range("AG1").ClearContents()

# --- Data row (row 2) fixes -----------------------------------------------
$ws.Range("A2").Value = "asd"
$ws.Range("AG2").Clear()

# --- Column width -----------------------------------------------------------
$ws.Columns.Item(16).ColumnWidth = 30

# --- Window scroll position + selection ------------------------------------
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 11
$ws.Range("P2").Select() | Out-Null
